# Auto-generated edit script implementing the "Week 12" journal entry addition
$d = $word.ActiveDocument

# 1. Remove the existing hidden _GoBack bookmark; we will recreate it spanning
#    the newly inserted content once that content exists.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$newContentStart = $d.Content.End

# --- paragraph 0 ---
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$p0 = $d.Paragraphs.Last.Range
$p0.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="PlainText"/>
        <w:spacing w:line="480" w:lineRule="auto"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Week 12 (April 24, 2016):</w:t>
      </w:r>
    </w:p>
'@) | Out-Null

# --- paragraph 1 ---
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$p1 = $d.Paragraphs.Last.Range
$p1.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="PlainText"/>
        <w:spacing w:line="480" w:lineRule="auto"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:tab/>
        <w:t>I started my final version of my bot. Worked on my DLL and created storage and timing inside of it. I have everything I need now to complete the project. The Windows Form application is as close as it gets to drag and drop, and managing tabs and text boxes is going to be extremely easy, although I wish it wasn’t as explicit.</w:t>
      </w:r>
    </w:p>
'@) | Out-Null

# --- paragraph 2 ---
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$p2 = $d.Paragraphs.Last.Range
$p2.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="PlainText"/>
        <w:spacing w:line="480" w:lineRule="auto"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:tab/>
        <w:t xml:space="preserve">Working with my own personal timer is going to be very irritating sometimes because Windows Forms doesn’t like it if you call other methods on other threads, as things can get messed up. This is where I learned that delegate methods and lambda expressions are things of </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>beauty. It’s like having an interface that you can carry around and put in as a parameter on everything, and just change and execute it whenever you please. This has been extremely helpful especially with my file editor and GUI manager.</w:t>
      </w:r>
    </w:p>
'@) | Out-Null

# --- paragraph 3 ---
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$p3 = $d.Paragraphs.Last.Range
$p3.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="PlainText"/>
        <w:spacing w:line="480" w:lineRule="auto"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:tab/>
        <w:t>Once everything gets rolling and I get the I/O situation figured out, everything else should be a breeze. Once I learn and pick out a language to use for storing my commands (I’m thinking JSON is the preferable choice), storage should be a very straightforward task.</w:t>
      </w:r>
    </w:p>
'@) | Out-Null

# --- paragraph 4 ---
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$p4 = $d.Paragraphs.Last.Range
$p4.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="PlainText"/>
        <w:spacing w:line="480" w:lineRule="auto"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:tab/>
        <w:t>Going through the final iteration has made me look back and realize how inexperienced I was back then with C#, even with Java, and that the best way to learn is to practice and struggle on your own. Luckily I have derived methods in which I don’t get errors due to my experience with stressing out (I’m looking at you parameter requiring user getting channel owner) for hours and toiling over a little mistake. Things get done quickly and I feel it’s not just me getting through the mud, but rather building tracks and starting the train. With the deliverable due next week, I will have my interactive bot that everyone so very much expects from me (I’m looking at you Isaac).</w:t>
      </w:r>
    </w:p>
'@) | Out-Null

# --- paragraph 5 ---
$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null
$p5 = $d.Paragraphs.Last.Range
$p5.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="PlainText"/>
        <w:spacing w:line="480" w:lineRule="auto"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>Total Word Count: 319</w:t>
      </w:r>
    </w:p>
'@) | Out-Null

# 2. Recreate the _GoBack bookmark spanning from just before the "Week 12"
#    heading run through the end of the newly added "Total Word Count: 319" run.
$bookmarkRange = $d.Range($newContentStart, $d.Content.End)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

Write-Host "Paragraphs now:" $d.Paragraphs.Count
